# roadtrip.xlsx - packing list / reservation updates (Option 1-Grand Canyon sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Option 1-Grand Canyon")
$ws.Activate() | Out-Null

# --- Content edits -------------------------------------------------------
# Row 7 used to describe the "key west / Bahia Honda" leg, row 8 the
# "leave Bahia Honda" leg; both are replaced by a short "pennekamp" note,
# and the adjacent Night column (E) now reads "reserved" for rows 6-8.
$ws.Range("D7").Value = "pennekamp"
$ws.Range("D8").Value = "pennekamp"
$ws.Range("E6").Value = "reserved"
$ws.Range("E7").Value = "reserved"
$ws.Range("E8").Value = "reserved"

# Rows 7 & 8 no longer need the taller wrapped-text row height now that the
# text is short -- let Excel recompute the (default) row height.
$ws.Rows.Item(7).AutoFit()
$ws.Rows.Item(8).AutoFit()

# --- View / selection state -----------------------------------------------
# Scroll the sheet so row 10 is at the top and select E6 (new active cell).
$win = $wb.Windows.Item(1)
$win.ScrollRow = 10
$win.ScrollColumn = 2
$ws.Range("E6").Select()

# Minimize the workbook window (best effort).
$win.WindowState = -4140
